# Update LR-pair TPM-derived values (Vegfa-Flt1) with newly computed TPM figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.436534333333333
$ws.Range("H2").Value = 4.309603
$ws.Range("I2").Value = 0.03241561610838976
$ws.Range("J2").Value = 0.03241561610838976
$ws.Range("M2").Value = 192.036265
$ws.Range("N2").Value = 576.108795
$ws.Range("O2").Value = 0.9537264487607444
$ws.Range("P2").Value = 0.9537264487607444
$ws.Range("Q2").Value = 275.8666879175983
$ws.Range("R2").Value = 2482.800191258385
$ws.Range("S2").Value = 0.03091563043544614
$ws.Range("T2").Value = 0.03091563043544614

# Row 3
$ws.Range("G3").Value = 1.436534333333333
$ws.Range("H3").Value = 4.309603
$ws.Range("I3").Value = 0.03241561610838976
$ws.Range("J3").Value = 0.03241561610838976
$ws.Range("O3").Value = 0.005072929450888834
$ws.Range("P3").Value = 0.005072929450888834
$ws.Range("Q3").Value = 1.467351825541556
$ws.Range("R3").Value = 13.206166429874
$ws.Range("S3").Value = 0.0001644421336249569
$ws.Range("T3").Value = 0.0001644421336249569

# Row 4
$ws.Range("G4").Value = 1.436534333333333
$ws.Range("H4").Value = 4.309603
$ws.Range("I4").Value = 0.03241561610838976
$ws.Range("J4").Value = 0.03241561610838976
$ws.Range("M4").Value = 8.295893999999999
$ws.Range("N4").Value = 24.887682
$ws.Range("O4").Value = 0.04120062178836673
$ws.Range("P4").Value = 0.04120062178836673
$ws.Range("Q4").Value = 11.917336556694
$ws.Range("R4").Value = 107.256029010246
$ws.Range("S4").Value = 0.001335543539318654
$ws.Range("T4").Value = 0.001335543539318654

# Row 5
$ws.Range("G5").Value = 34.88211266666666
$ws.Range("I5").Value = 0.7871201871162607
$ws.Range("J5").Value = 0.7871201871162609
$ws.Range("M5").Value = 192.036265
$ws.Range("N5").Value = 576.108795
$ws.Range("O5").Value = 0.9537264487607444
$ws.Range("P5").Value = 0.9537264487607444
$ws.Range("Q5").Value = 6698.630631815856
$ws.Range("R5").Value = 60287.67568634271
$ws.Range("S5").Value = 0.750697340806284
$ws.Range("T5").Value = 0.7506973408062841

# Row 6
$ws.Range("G6").Value = 34.88211266666666
$ws.Range("I6").Value = 0.7871201871162607
$ws.Range("J6").Value = 0.7871201871162609
$ws.Range("O6").Value = 0.005072929450888834
$ws.Range("P6").Value = 0.005072929450888834
$ws.Range("Q6").Value = 35.63042700233378
$ws.Range("R6").Value = 320.673843021004
$ws.Range("S6").Value = 0.003993005178611208
$ws.Range("T6").Value = 0.003993005178611209

# Row 7
$ws.Range("G7").Value = 34.88211266666666
$ws.Range("I7").Value = 0.7871201871162607
$ws.Range("J7").Value = 0.7871201871162609
$ws.Range("M7").Value = 8.295893999999999
$ws.Range("N7").Value = 24.887682
$ws.Range("O7").Value = 0.04120062178836673
$ws.Range("P7").Value = 0.04120062178836673
$ws.Range("Q7").Value = 289.378309178724
$ws.Range("R7").Value = 2604.404782608516
$ws.Range("S7").Value = 0.03242984113136551
$ws.Range("T7").Value = 0.03242984113136551

# Row 8
$ws.Range("G8").Value = 7.997472999999999
$ws.Range("H8").Value = 23.992419
$ws.Range("I8").Value = 0.1804641967753495
$ws.Range("J8").Value = 0.1804641967753495
$ws.Range("M8").Value = 192.036265
$ws.Range("N8").Value = 576.108795
$ws.Range("O8").Value = 0.9537264487607444
$ws.Range("P8").Value = 0.9537264487607444
$ws.Range("Q8").Value = 1535.804844358345
$ws.Range("R8").Value = 13822.2435992251
$ws.Range("S8").Value = 0.1721134775190142
$ws.Range("T8").Value = 0.1721134775190142

# Row 9
$ws.Range("G9").Value = 7.997472999999999
$ws.Range("H9").Value = 23.992419
$ws.Range("I9").Value = 0.1804641967753495
$ws.Range("J9").Value = 0.1804641967753495
$ws.Range("O9").Value = 0.005072929450888834
$ws.Range("P9").Value = 0.005072929450888834
$ws.Range("Q9").Value = 8.169040122444667
$ws.Range("R9").Value = 73.521361102002
$ws.Range("S9").Value = 0.000915482138652668
$ws.Range("T9").Value = 0.000915482138652668

# Row 10
$ws.Range("G10").Value = 7.997472999999999
$ws.Range("H10").Value = 23.992419
$ws.Range("I10").Value = 0.1804641967753495
$ws.Range("J10").Value = 0.1804641967753495
$ws.Range("M10").Value = 8.295893999999999
$ws.Range("N10").Value = 24.887682
$ws.Range("O10").Value = 0.04120062178836673
$ws.Range("P10").Value = 0.04120062178836673
$ws.Range("Q10").Value = 66.34618827586199
$ws.Range("R10").Value = 597.1156944827579
$ws.Range("S10").Value = 0.007435237117682563
$ws.Range("T10").Value = 0.007435237117682563
